# ExcelRow.copyAndInsert: shift only if not at the end.
# ExcelSheet: getEmptyRows < lastRowNum.
#
# Workbook-Test.xlsx / "Row-actions" sheet: renumber the row-index column
# (A) to reflect that rows are only shifted down when the insertion point
# is not already at the end of the sheet, and drop the now-unused external
# link to Test.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Row-actions")

$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8

$ws.Activate()
$ws.Range("A9").Select()

# Remove the stale external reference to Test.xlsx that is no longer needed.
$links = $wb.LinkSources()
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}
